# Update "想去人数" (interest count) values in column F for both the
# "展览" (sheet 1) and "全部类型" (sheet 4) worksheets, which hold the
# same underlying data.

$wb = $excel.ActiveWorkbook

# Map of cell address -> new value to apply on each target worksheet.
$updates = @{
    "F2"  = 2019
    "F7"  = 1681
    "F9"  = 679
    "F12" = 26
    "F13" = 98
    "F14" = 227
    "F19" = 3898
    "F22" = 438
    "F23" = 363
    "F24" = 726
    "F25" = 548
    "F26" = 357
    "F27" = 34
    "F28" = 1693
    "F29" = 16
    "F30" = 27
    "F31" = 166
}

$targetSheetIndexes = @(1, 4)

foreach ($sheetIndex in $targetSheetIndexes) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
